# Apply "Add data for 2022-08-23" update to the carjacking-by-neighborhood-by-month workbook.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Rename the sheet (tab name) and update the header label cell (B1) to reflect
# the new "through" date.
$ws.Name = "Through 2022-08-15"
$ws.Range("B1").Value = "August 2022 (through August 15)"

# Update existing cell values (counts incremented because of newly-added data).
$ws.Range("R2").Value = 4
$ws.Range("R3").Value = 3
$ws.Range("D4").Value = 6
$ws.Range("Z5").Value = 2
$ws.Range("AX5").Value = 2
$ws.Range("J6").Value = 3
$ws.Range("B8").Value = 2
$ws.Range("B9").Value = 2
$ws.Range("B13").Value = 2
$ws.Range("R13").Value = 2
$ws.Range("AP17").Value = 3
$ws.Range("AP20").Value = 2
$ws.Range("R28").Value = 3
$ws.Range("B29").Value = 6
$ws.Range("J29").Value = 4
$ws.Range("B37").Value = 3
$ws.Range("J39").Value = 2
$ws.Range("AX46").Value = 2
$ws.Range("AX66").Value = 2

# Fill in cells that previously had no value (new data points).
$ws.Range("R4").Value = 1
$ws.Range("AH14").Value = 1
$ws.Range("R34").Value = 1
$ws.Range("J50").Value = 2
$ws.Range("R74").Value = 1
$ws.Range("J75").Value = 1
$ws.Range("J90").Value = 1
$ws.Range("BF92").Value = 1
